$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value2 = 361.33334
$ws.Range("I12").Value2 = 380.54544
$ws.Range("J12").Value2 = 150
$ws.Range("K12").Value2 = 380.54544
$ws.Range("L12").Value2 = 150
$ws.Range("M12").Value2 = -210.54544
$ws.Range("N12").Value2 = -490

$ws.Range("H15").Value2 = 1450.1923
$ws.Range("I15").Value2 = 1450.1923
$ws.Range("K15").Value2 = 4350.5769
$ws.Range("M15").Value2 = -4181.5769

$ws.Range("H33").Value2 = 465.29413
$ws.Range("I33").Value2 = 294.06668
$ws.Range("K33").Value2 = 294.06668
$ws.Range("M33").Value2 = -65.06668000000002

$ws.Range("H48").Value2 = 12112.6
$ws.Range("J48").Value2 = 13854.333
$ws.Range("L48").Value2 = 41562.999
$ws.Range("N48").Value2 = -42146.999

$ws.Range("H51").Value2 = 5842.0527
$ws.Range("I51").Value2 = 6249.9375
$ws.Range("J51").Value2 = 3666.6667
$ws.Range("K51").Value2 = 6249.9375
$ws.Range("L51").Value2 = 3666.6667
$ws.Range("M51").Value2 = -5765.9375
$ws.Range("N51").Value2 = -4634.6667

$ws.Range("H53").Value2 = 1166.2307
$ws.Range("I53").Value2 = 828.8570999999999
$ws.Range("J53").Value2 = 1559.8334
$ws.Range("K53").Value2 = 828.8570999999999
$ws.Range("L53").Value2 = 1559.8334
$ws.Range("M53").Value2 = -191.8570999999999
$ws.Range("N53").Value2 = -2833.8334

$ws.Range("H56").Value2 = 12112.6
$ws.Range("J56").Value2 = 13854.333
$ws.Range("L56").Value2 = 41562.999
$ws.Range("N56").Value2 = -42630.999

$ws.Range("H113").Value2 = 4437.5815
$ws.Range("I113").Value2 = 5112.697
$ws.Range("J113").Value2 = 2209.7
$ws.Range("K113").Value2 = 5112.697
$ws.Range("L113").Value2 = 2209.7
$ws.Range("M113").Value2 = -1858.697
$ws.Range("N113").Value2 = -8717.700000000001

$ws.Range("H132").Value2 = 1926.9756
$ws.Range("I132").Value2 = 2000.6052
$ws.Range("K132").Value2 = 6001.8156
$ws.Range("M132").Value2 = -3471.8156

$ws.Range("H141").Value2 = 0
$ws.Range("I141").Value2 = 0
$ws.Range("K141").Value2 = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value2 = 1711.4615
$ws.Range("I45").Value2 = 1113.5454
$ws.Range("J45").Value2 = 5000
$ws.Range("K45").Value2 = 1113.5454
$ws.Range("L45").Value2 = 5000
$ws.Range("M45").Value2 = -736.5454
$ws.Range("N45").Value2 = -5754

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H116").Value2 = 61495
$ws.Range("J116").Value2 = 61495
$ws.Range("L116").Value2 = 61495
$ws.Range("N116").Value2 = -70673

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value2 = 504.85715
$ws.Range("I22").Value2 = 504.85715
$ws.Range("J22").Value2 = 0
$ws.Range("K22").Value2 = 504.85715
$ws.Range("L22").Value2 = 0
$ws.Range("M22").Value2 = -154.85715
$ws.Range("N22").ClearContents()

$ws.Range("H31").Value2 = 6689.184
$ws.Range("I31").Value2 = 2609.4614
$ws.Range("J31").Value2 = 8810.639999999999
$ws.Range("K31").Value2 = 2609.4614
$ws.Range("L31").Value2 = 8810.639999999999
$ws.Range("M31").Value2 = -2314.4614
$ws.Range("N31").Value2 = -9400.639999999999

$ws.Range("H34").Value2 = 6689.184
$ws.Range("I34").Value2 = 2609.4614
$ws.Range("J34").Value2 = 8810.639999999999
$ws.Range("K34").Value2 = 2609.4614
$ws.Range("L34").Value2 = 8810.639999999999
$ws.Range("M34").Value2 = -2407.4614
$ws.Range("N34").Value2 = -9214.639999999999

$ws.Range("H111").Value2 = 120000
$ws.Range("J111").Value2 = 120000
$ws.Range("L111").Value2 = 120000
$ws.Range("N111").Value2 = -128180

$ws.Range("H132").Value2 = 3178.9092
$ws.Range("I132").Value2 = 2890.7334
$ws.Range("K132").Value2 = 8672.200199999999
$ws.Range("M132").Value2 = -6142.200199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value2 = 1699.7333
$ws.Range("I5").Value2 = 1148.8182
$ws.Range("K5").Value2 = 3446.4546
$ws.Range("M5").Value2 = -3334.4546

$ws.Range("H75").Value2 = 16165.333
$ws.Range("J75").Value2 = 18598.4
$ws.Range("L75").Value2 = 55795.2
$ws.Range("N75").Value2 = -57791.2

$ws.Range("H78").Value2 = 16165.333
$ws.Range("J78").Value2 = 18598.4
$ws.Range("L78").Value2 = 167385.6
$ws.Range("N78").Value2 = -177369.6

$ws.Range("H117").Value2 = 4216.273
$ws.Range("J117").Value2 = 4753
$ws.Range("L117").Value2 = 14259
$ws.Range("N117").Value2 = -21143

$ws.Range("H121").Value2 = 39808.555
$ws.Range("J121").Value2 = 44676.5
$ws.Range("L121").Value2 = 134029.5
$ws.Range("N121").Value2 = -136649.5

$ws.Range("H135").Value2 = 1699.7333
$ws.Range("I135").Value2 = 1148.8182
$ws.Range("K135").Value2 = 10339.3638
$ws.Range("M135").Value2 = -7804.363799999999

$ws.Range("H137").Value2 = 12458.917
$ws.Range("I137").Value2 = 1277.5
$ws.Range("J137").Value2 = 18049.625
$ws.Range("K137").Value2 = 3832.5
$ws.Range("L137").Value2 = 54148.875
$ws.Range("M137").Value2 = 1267.5
$ws.Range("N137").Value2 = -64348.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value2 = 4245034.5
$ws.Range("J11").Value2 = 1790903.4
$ws.Range("L11").Value2 = 1790903.4
$ws.Range("N11").Value2 = -1791181.4

$ws.Range("H62").Value2 = 51996.668
$ws.Range("J62").Value2 = 55995
$ws.Range("L62").Value2 = 55995
$ws.Range("N62").Value2 = -57367

$ws.Range("H65").Value2 = 51996.668
$ws.Range("J65").Value2 = 55995
$ws.Range("L65").Value2 = 167985
$ws.Range("N65").Value2 = -174849

$ws.Range("H80").Value2 = 3007.375
$ws.Range("I80").Value2 = 2576
$ws.Range("J80").Value2 = 3438.75
$ws.Range("K80").Value2 = 2576
$ws.Range("L80").Value2 = 3438.75
$ws.Range("M80").Value2 = -1578
$ws.Range("N80").Value2 = -5434.75

$ws.Range("H83").Value2 = 3007.375
$ws.Range("I83").Value2 = 2576
$ws.Range("J83").Value2 = 3438.75
$ws.Range("K83").Value2 = 12880
$ws.Range("L83").Value2 = 17193.75
$ws.Range("M83").Value2 = -7888
$ws.Range("N83").Value2 = -27177.75

$ws.Range("H113").Value2 = 16171.714
$ws.Range("I113").Value2 = 2095
$ws.Range("J113").Value2 = 34940.668
$ws.Range("K113").Value2 = 2095
$ws.Range("L113").Value2 = 34940.668
$ws.Range("M113").Value2 = 75
$ws.Range("N113").Value2 = -39280.668

$ws.Range("H116").Value2 = 98257.664
$ws.Range("J116").Value2 = 98257.664
$ws.Range("L116").Value2 = 98257.664
$ws.Range("N116").Value2 = -107435.664

$ws.Range("H124").Value2 = 122997.5
$ws.Range("J124").Value2 = 122997.5
$ws.Range("L124").Value2 = 122997.5
$ws.Range("N124").Value2 = -132817.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value2 = 300
$ws.Range("I16").Value2 = 300
$ws.Range("J16").Value2 = 0
$ws.Range("K16").Value2 = 300
$ws.Range("L16").Value2 = 0
$ws.Range("M16").Value2 = -130
$ws.Range("N16").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value2 = 16839.4
$ws.Range("J51").Value2 = 20999
$ws.Range("L51").Value2 = 20999
$ws.Range("N51").Value2 = -22019

$ws.Range("H132").Value2 = 1284.8108
$ws.Range("I132").Value2 = 1347.4839
$ws.Range("J132").Value2 = 961
$ws.Range("K132").Value2 = 4042.4517
$ws.Range("L132").Value2 = 2883
$ws.Range("M132").Value2 = -1512.4517
$ws.Range("N132").Value2 = -7943

$ws.Range("H136").Value2 = 1350.258
$ws.Range("I136").Value2 = 1075.3214
$ws.Range("J136").Value2 = 3916.3333
$ws.Range("K136").Value2 = 3916.3333
$ws.Range("L136").Value2 = 11748.9999
$ws.Range("M136").Value2 = -675.9642000000003
$ws.Range("N136").Value2 = -16848.9999
